$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colmap = @{ A=1; B=2; C=3; D=4; E=5; F=6; G=7; H=8; I=9; J=10; K=11; L=12; M=13; N=14; O=15; P=16; Q=17; R=18; S=19; T=20; U=21; V=22; W=23; X=24; Y=25 }

$rows = @{}

$rows[82] = @{ A="November08  19:23:17"; B="resnet18"; C="'200"; D="'1"; E="'True"; F="sgd"; G="1.0xsingle + 0.0Xmulti"; H="'0.01"; I="'0.9"; J="<function exp_lr_scheduler at 0x7f07ec316ae8>"; K="'10"; L="'True"; M="'32"; S=" " }
$rows[83] = @{ A="November08  19:23:19"; B="resnet18"; C="'200"; D="'1"; E="'True"; F="sgd"; G="1.0xsingle + 0.0Xmulti"; H="'0.01"; I="'0.9"; J="<function exp_lr_scheduler at 0x7f07ec316ae8>"; K="'10"; L="'True"; M="'32"; N=5; O=0.04334113671348429; P=0.0281702742790267; Q=0.5352798053527981; R=0.6747769667477697; S=1.388751006943395; T=1.111030005150121; U=5; V=0.5352798053527981; W=0.6747769667477697; X=1.388751006943395; Y=1.111030005150121 }
$rows[84] = @{ A="November08  19:26:51"; B="resnet18"; C="'200"; D="'1"; E="'True"; F="sgd"; G="1.0xsingle + 0.0Xmulti"; H="'0.01"; I="'0.9"; J="<function exp_lr_scheduler at 0x7f07ec316ae8>"; K="'10"; L="'True"; M="'32"; N=2; O=0.0419952683785246; P=0.01750046590045518; Q=0.3098134630981346; R=0.07137064071370641; S=1.379081349673997; T=0.8054961244402359; U=2; V=0.3098134630981346; W=0.07137064071370641; X=1.379081349673997; Y=0.8054961244402359 }
$rows[85] = @{ A="November08  19:29:56"; B="resnet18"; C="'200"; D="'1"; E="'True"; F="sgd"; G="1.0xsingle + 0.0Xmulti"; H="'0.01"; I="'0.9"; J="<function exp_lr_scheduler at 0x7f07ec316ae8>"; K="'10"; L="'True"; M="'32"; N=20; O=0.002690032043885359; P=0.07745679709818457; Q=0.9816017316017316; R=0.4090909090909091; S=0.2326210525996177; T=1.421084275995101; U=20; V=0.9816017316017316; W=0.4090909090909091; X=0.2326210525996177; Y=1.421084275995101 }
$rows[86] = @{ A="November08  19:35:46"; B="resnet18"; C="'200"; D="'1"; E="'True"; F="sgd"; G="1.0xsingle + 0.0Xmulti"; H="'0.01"; I="'0.9"; J="<function exp_lr_scheduler at 0x7f07ec316ae8>"; K="'10"; L="'True"; M="'32"; N=26; O=0.00232736115622056; P=0.07711041399410792; Q=0.9816017316017316; R=0.3863636363636364; S=0.3068482428756015; T=1.403844443226963; U=26; V=0.9816017316017316; W=0.3863636363636364; X=0.3068482428756015; Y=1.403844443226963 }
$rows[87] = @{ A="November08  19:41:39"; B="resnet18"; C="'200"; D="'1"; E="'True"; F="sgd"; G="1.0xsingle + 0.0Xmulti"; H="'0.01"; I="'0.9"; J="<function exp_lr_scheduler at 0x7f07ec316ae8>"; K="'10"; L="'True"; M="'32"; N=24; O=0.003207404706707764; P=0.08552800796248695; Q=0.9761904761904762; R=0.3766233766233766; S=0.3512500866571044; T=1.536355951721349; U=24; V=0.9761904761904762; W=0.3766233766233766; X=0.3512500866571044; Y=1.536355951721349 }
$rows[88] = @{ A="November08  19:47:31"; B="resnet18"; C="'200"; D="'1"; E="'True"; F="sgd"; G="1.0xsingle + 0.0Xmulti"; H="'0.01"; I="'0.9"; J="<function exp_lr_scheduler at 0x7f07ec316ae8>"; K="'10"; L="'True"; M="'32"; N=22; O=0.002484592242222844; P=0.07642197802469328; Q=0.974025974025974; R=0.3993506493506493; S=0.3434608681718603; T=1.456066574758816; U=22; V=0.974025974025974; W=0.3993506493506493; X=0.3434608681718603; Y=1.456066574758816 }
$rows[89] = @{ A="November08  19:53:23"; B="resnet18"; C="'200"; D="'1"; E="'True"; F="sgd"; G="0.0xsingle + 1.0Xmulti"; H="'0.01"; I="'0.9"; J="<function exp_lr_scheduler at 0x7f07ec316ae8>"; K="'10"; L="'True"; M="'32"; N=14; O=0.00454240246349341; P=0.01250092022411235; Q=0.5757575757575758; R=0.3636363636363636; S=0.719246840323972; T=1.237928726731686; U=14; V=0.5757575757575758; W=0.3636363636363636; X=0.719246840323972; Y=1.237928726731686 }
$rows[90] = @{ A="November08  19:59:12"; B="resnet18"; C="'200"; D="'1"; E="'True"; F="sgd"; G="0.0xsingle + 1.0Xmulti"; H="'0.01"; I="'0.9"; J="<function exp_lr_scheduler at 0x7f07ec316ae8>"; K="'10"; L="'True"; M="'32"; N=7; O=0.007189551182897576; P=0.01224703931963289; Q=0.4675324675324675; R=0.3311688311688312; S=0.8765832923335172; T=1.289316742440608; U=7; V=0.4675324675324675; W=0.3311688311688312; X=0.8765832923335172; Y=1.289316742440608 }
$rows[91] = @{ A="November08  20:05:03"; B="resnet18"; C="'200"; D="'1"; E="'True"; F="sgd"; G="0.0xsingle + 1.0Xmulti"; H="'0.01"; I="'0.9"; J="<function exp_lr_scheduler at 0x7f07ec316ae8>"; K="'10"; L="'True"; M="'32"; N=19; O=0.003655296125453272; P=0.01348998949125216; Q=0.6006493506493507; R=0.3636363636363636; S=0.6546536707079771; T=1.310546629285836; U=19; V=0.6006493506493507; W=0.3636363636363636; X=0.6546536707079771; Y=1.310546629285836 }
$rows[92] = @{ A="November08  20:10:56"; B="resnet18"; C="'200"; D="'1"; E="'True"; F="sgd"; G="0.0xsingle + 1.0Xmulti"; H="'0.01"; I="'0.9"; J="<function exp_lr_scheduler at 0x7f07ec316ae8>"; K="'10"; L="'True"; M="'32"; N=16; O=0.00432416656381124; P=0.01360972293398597; Q=0.5573593073593074; R=0.3733766233766234; S=0.7222499717168716; T=1.258735708723473; U=16; V=0.5573593073593074; W=0.3733766233766234; X=0.7222499717168716; Y=1.258735708723473 }
$rows[93] = @{ A="November08  20:16:48"; B="resnet18"; C="'200"; D="'1"; E="'True"; F="sgd"; G="0.1xsingle + 0.9Xmulti"; H="'0.01"; I="'0.9"; J="<function exp_lr_scheduler at 0x7f07ec316ae8>"; K="'10"; L="'True"; M="'32"; N=19; O=0.004559562214976782; P=0.01790383312996332; Q=0.8885281385281385; R=0.4318181818181818; S=0.5243012253388074; T=1.306825239958353; U=19; V=0.8885281385281385; W=0.4318181818181818; X=0.5243012253388074; Y=1.306825239958353 }
$rows[94] = @{ A="November08  20:22:45"; B="resnet18"; C="'200"; D="'1"; E="'True"; F="sgd"; G="0.1xsingle + 0.9Xmulti"; H="'0.01"; I="'0.9"; J="<function exp_lr_scheduler at 0x7f07ec316ae8>"; K="'10"; L="'True"; M="'32"; N=30; O=0.0036955579341232; P=0.01761063381836012; Q=0.908008658008658; R=0.4448051948051948; S=0.3678061789603123; T=1.1870513506546; U=30; V=0.908008658008658; W=0.4448051948051948; X=0.3678061789603123; Y=1.1870513506546 }
$rows[95] = @{ A="November08  20:28:38"; B="resnet18"; C="'200"; D="'1"; E="'True"; F="sgd"; G="0.1xsingle + 0.9Xmulti"; H="'0.01"; I="'0.9"; J="<function exp_lr_scheduler at 0x7f07ec316ae8>"; K="'10"; L="'True"; M="'32"; N=14; O=0.004371386480989394; P=0.01904632364000593; Q=0.8766233766233766; R=0.4318181818181818; S=0.4912657034049816; T=1.28427047359606; U=14; V=0.8766233766233766; W=0.4318181818181818; X=0.4912657034049816; Y=1.28427047359606 }
$rows[96] = @{ A="November08  20:34:32"; B="resnet18"; C="'200"; D="'1"; E="'True"; F="sgd"; G="0.1xsingle + 0.9Xmulti"; H="'0.01"; I="'0.9"; J="<function exp_lr_scheduler at 0x7f07ec316ae8>"; K="'10"; L="'True"; M="'32"; N=25; O=0.003952494100768329; P=0.02005540086077405; Q=0.8852813852813853; R=0.3863636363636364; S=0.4301917422091399; T=1.245772070624499; U=25; V=0.8852813852813853; W=0.3863636363636364; X=0.4301917422091399; Y=1.245772070624499 }

foreach ($rowNum in $rows.Keys) {
    $rowData = $rows[$rowNum]
    foreach ($key in $rowData.Keys) {
        $col = $colmap[$key]
        $ws.Cells.Item([int]$rowNum, $col).Value = $rowData[$key]
    }
}

Write-Host "done writing rows"